$wb = $excel.ActiveWorkbook

$wsTastes = $wb.Worksheets.Item("Tastes & aromas of varieties")
$wsFood   = $wb.Worksheets.Item("Wine & food pairing")

# --- Data entry, replicated in the same order the shared-string table was
# --- built (first-use order), so new strings land at the expected indices.

# Row 1 of the food-pairing sheet (Vinho Verde Albariño)
$wsFood.Range("B1").Value = "Clams, Octopus Salad"
$wsFood.Range("A1").Value = "Vinho Verde Albariño"

# Assyrtiko
$wsTastes.Range("A6").Value = "Assyrtiko"
$wsTastes.Range("A6").WrapText = $true
$wsTastes.Range("B6").Value = "Good acidity"
$wsFood.Range("B2").Value = "Seafood, fishes"

# Fiano di Avellino
$wsTastes.Range("A7").Value = "Fiano di Avellino"
$wsTastes.Range("B7").Value = "Delicate, subtil, with flower flavours"
$wsFood.Range("B3").Value = "Fish, and seafood pastas"

# Greco di Tufo
$wsTastes.Range("A8").Value = "Greco di Tufo"
$wsTastes.Range("B8").Value = "Very fresh, mineral, intense, white flowers and white fruits aromas"
$wsFood.Range("B4").Value = "Shellfishes, Fish with white flesh"

# Gruner Veltliner
$wsTastes.Range("A9").Value = "Grüner Veltliner"
$wsTastes.Range("B9").Value = "White pepper, mellifère (nectar plants) aromas"
$wsFood.Range("B5").Value = "White meat, wild game meat"

# Petite arvine
$wsTastes.Range("A10").Value = "Petite arvine"
$wsTastes.Range("B10").Value = "Exotic or salin aromas"
$wsFood.Range("B6").Value = "River fishes, Exotic fruits soufflés, Safran ice cream, White truffle risotto"

# Stray numeric value further down the "Tastes & aromas of varieties" sheet
$wsTastes.Range("A15").Value = 1410

# Fill in the wine-variety names on the food-pairing sheet (reuses strings
# already present in the shared-string table, so no new entries here)
$wsFood.Range("A2").Value = "Assyrtiko"
$wsFood.Range("A2").WrapText = $true
$wsFood.Range("A3").Value = "Fiano di Avellino"
$wsFood.Range("A4").Value = "Greco di Tufo"
$wsFood.Range("A5").Value = "Grüner Veltliner"
$wsFood.Range("A6").Value = "Petite arvine"

# --- Leave the edit cursor where the author left it on each sheet
$wsTastes.Range("A16").Select() | Out-Null
$wsFood.Range("A11").Select() | Out-Null

# --- Re-order the tabs: "Wine & food pairing" now sits before
# --- "Varieties of designations"
$wsVarieties = $wb.Worksheets.Item("Varieties of designations")
$wsFood.Move($wsVarieties)

# --- Re-activate "Tastes & aromas of varieties" (was, and stays, the
# --- selected tab)
$wsTastes.Activate()
